$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Duplicate Sheet1 ("Move or Copy" -> "Create a copy") to make Sheet2,
# placed immediately after Sheet1.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Turn on AutoFilter for Sheet2's data range.
$ws2.Range("A1:D49").AutoFilter()

# Register the (hidden) _FilterDatabase defined name that Excel creates
# for a filtered range, scoped to Sheet2.
$fdb = $ws2.Names.Add("_xlnm._FilterDatabase", "=Sheet2!`$A`$1:`$D`$49")
$fdb.Visible = $false

# Restore the on-screen selections: Sheet2 left with G12 selected (and not
# the active/tab-selected sheet), Sheet1 remains the active sheet with B15
# selected.
$ws2.Range("G12").Select()
$ws1.Activate()
$ws1.Range("B15").Select()
